$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text format before writing, so values like
# "582.72" or "171.17" are stored as text (matching the original inline-string
# cells) instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '66.508.55'
$ws.Range("E2").Value = '  -1.75%  '
$ws.Range("D3").Value = '2.509.07'
$ws.Range("E3").Value = '  -4.85%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '582.72'
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").Value = '171.17'
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '0.525'
$ws.Range("E8").Value = '  -1.84%  '
$ws.Range("D9").Value = '2.509.22'
$ws.Range("E9").Value = '  -4.88%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '0.350'
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("D13").Value = '5.12'
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").Value = '26.68'
$ws.Range("E14").Value = '  -3.54%  '
$ws.Range("D15").Value = '2.957.24'
$ws.Range("E15").Value = '  -4.96%  '
$ws.Range("E16").Value = '  -3.36%  '
$ws.Range("D17").Value = '66.370.28'
$ws.Range("E17").Value = '  -4.19%  '
$ws.Range("D18").Value = '2.501.89'
$ws.Range("E18").Value = '  -4.06%  '
$ws.Range("D19").Value = '11.24'
$ws.Range("E19").Value = '  -6.86%  '
$ws.Range("D20").Value = '7.70'
$ws.Range("E20").Value = '  -4.75%  '
$ws.Range("D21").Value = '347.23'
$ws.Range("E21").Value = '  -3.00%  '
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  -3.31%  '
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("D24").Value = '1.97'
$ws.Range("E24").Value = '  +1.06%  '
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").Value = '69.60'
$ws.Range("E26").Value = '  -0.63%  '
$ws.Range("D27").Value = '9.95'
$ws.Range("E27").Value = '  -4.24%  '
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("E29").Value = '  -4.91%  '
$ws.Range("D30").Value = '0.0₃0976'
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").Value = '526.19'
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("D32").Value = '8.11'
$ws.Range("E32").Value = '  +1.70%  '
$ws.Range("E33").Value = '  -3.28%  '
$ws.Range("E34").Value = '  -3.65%  '
$ws.Range("D35").Value = '0.131'
$ws.Range("E35").Value = '  -4.66%  '
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '1.46'
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").Value = '155.50'
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").Value = '18.60'
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("D40").Value = '18.37'
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("E41").Value = '  -3.40%  '
$ws.Range("D42").Value = '1.79'
$ws.Range("E42").Value = '  -1.89%  '
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("E44").Value = '  -0.10%  '
$ws.Range("E45").Value = '  +2.45%  '
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").Value = '148.33'
$ws.Range("E47").Value = '  -3.16%  '
$ws.Range("D48").Value = '0.558'
$ws.Range("E48").Value = '  -4.22%  '
$ws.Range("D49").Value = '3.68'
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("D50").Value = '0.0₆0272'
$ws.Range("E50").Value = '  -9.49%  '
$ws.Range("D51").Value = '1.72'
$ws.Range("E51").Value = '  +0.34%  '

# Restore the default (Normal) cell style on the Price column so the
# unaffected style index bookkeeping matches a normal Excel edit and cells
# fall back to the workbook's default (unstyled) formatting.
$ws.Range("D2:D51").Style = "Normal"
